$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 5182.8945
$ws.Range("I100").Value = 2672.8572
$ws.Range("J100").Value = 6647.0835
$ws.Range("K100").Value = 2672.8572
$ws.Range("L100").Value = 6647.0835
$ws.Range("M100").Value = -2131.8572
$ws.Range("N100").Value = -7729.0835
$ws.Range("H116").Value = 2188
$ws.Range("I116").Value = 2208.3
$ws.Range("J116").Value = 1985
$ws.Range("K116").Value = 2208.3
$ws.Range("L116").Value = 1985
$ws.Range("M116").Value = 1233.7
$ws.Range("N116").Value = -8869
$ws.Range("H118").Value = 2079.8
$ws.Range("I118").Value = 2366.5
$ws.Range("J118").Value = 1649.75
$ws.Range("K118").Value = 7099.5
$ws.Range("L118").Value = 4949.25
$ws.Range("M118").Value = -5442.5
$ws.Range("N118").Value = -8263.25
$ws.Range("H132").Value = 16390.576
$ws.Range("I132").Value = 1175
$ws.Range("K132").Value = 3525
$ws.Range("M132").Value = -995
$ws.Range("H137").Value = 1891.3704
$ws.Range("I137").Value = 1281.2632
$ws.Range("K137").Value = 3843.7896
$ws.Range("M137").Value = -1293.7896

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H39").Value = 1000
$ws.Range("I39").Value = 1000
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 1000
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -480
$ws.Range("N39").ClearContents()
$ws.Range("H74").Value = 2559.5454
$ws.Range("I74").Value = 2193.4119
$ws.Range("J74").Value = 3804.4
$ws.Range("K74").Value = 2193.4119
$ws.Range("L74").Value = 3804.4
$ws.Range("M74").Value = -1319.4119
$ws.Range("N74").Value = -5552.4
$ws.Range("H77").Value = 2559.5454
$ws.Range("I77").Value = 2193.4119
$ws.Range("J77").Value = 3804.4
$ws.Range("K77").Value = 10967.0595
$ws.Range("L77").Value = 19022
$ws.Range("M77").Value = -6599.059499999999
$ws.Range("N77").Value = -27758
$ws.Range("H88").Value = 1894
$ws.Range("I88").Value = 1796.3334
$ws.Range("J88").Value = 1967.25
$ws.Range("K88").Value = 1796.3334
$ws.Range("L88").Value = 1967.25
$ws.Range("M88").Value = -1390.3334
$ws.Range("N88").Value = -2779.25
$ws.Range("H91").Value = 1894
$ws.Range("I91").Value = 1796.3334
$ws.Range("J91").Value = 1967.25
$ws.Range("K91").Value = 1796.3334
$ws.Range("L91").Value = 1967.25
$ws.Range("M91").Value = -392.3334
$ws.Range("N91").Value = -4775.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 76497.91
$ws.Range("J138").Value = 76560.48
$ws.Range("L138").Value = 76560.48
$ws.Range("N138").Value = -86840.48

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1459.25
$ws.Range("J16").Value = 2497.5
$ws.Range("L16").Value = 2497.5
$ws.Range("N16").Value = -3071.5
$ws.Range("H31").Value = 5034.4116
$ws.Range("I31").Value = 2485.5715
$ws.Range("J31").Value = 6818.6
$ws.Range("K31").Value = 2485.5715
$ws.Range("L31").Value = 6818.6
$ws.Range("M31").Value = -2190.5715
$ws.Range("N31").Value = -7408.6
$ws.Range("H34").Value = 5034.4116
$ws.Range("I34").Value = 2485.5715
$ws.Range("J34").Value = 6818.6
$ws.Range("K34").Value = 2485.5715
$ws.Range("L34").Value = 6818.6
$ws.Range("M34").Value = -2283.5715
$ws.Range("N34").Value = -7222.6
$ws.Range("H113").Value = 1459.25
$ws.Range("J113").Value = 2497.5
$ws.Range("L113").Value = 2497.5
$ws.Range("N113").Value = -6837.5
$ws.Range("H122").Value = 45457390
$ws.Range("J122").Value = 4686
$ws.Range("L122").Value = 14058
$ws.Range("N122").Value = -18958

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 600.43475
$ws.Range("I114").Value = 290.54544
$ws.Range("K114").Value = 871.63632
$ws.Range("M114").Value = 2382.36368

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 1000
$ws.Range("I32").Value = 1000
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1000
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -683
$ws.Range("N32").ClearContents()
$ws.Range("H59").Value = 35000
$ws.Range("J59").Value = 35000
$ws.Range("L59").Value = 35000
$ws.Range("N59").Value = -36308
$ws.Range("H61").Value = 1846.1892
$ws.Range("I61").Value = 1244.9546
$ws.Range("J61").Value = 2728
$ws.Range("K61").Value = 1244.9546
$ws.Range("L61").Value = 2728
$ws.Range("M61").Value = -1042.9546
$ws.Range("N61").Value = -3132
$ws.Range("H113").Value = 1846.1892
$ws.Range("I113").Value = 1244.9546
$ws.Range("J113").Value = 2728
$ws.Range("K113").Value = 1244.9546
$ws.Range("L113").Value = 2728
$ws.Range("M113").Value = 925.0454
$ws.Range("N113").Value = -7068
$ws.Range("H136").Value = 7139.579
$ws.Range("I136").Value = 5986
$ws.Range("J136").Value = 7812.5
$ws.Range("K136").Value = 17958
$ws.Range("L136").Value = 23437.5
$ws.Range("M136").Value = -15408
$ws.Range("N136").Value = -28537.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 10500
$ws.Range("I21").Value = 10500
$ws.Range("K21").Value = 10500
$ws.Range("M21").Value = -10265
$ws.Range("H33").Value = 12331.667
$ws.Range("I33").Value = 3750
$ws.Range("K33").Value = 3750
$ws.Range("M33").Value = -3500
$ws.Range("H35").Value = 10500
$ws.Range("I35").Value = 10500
$ws.Range("K35").Value = 10500
$ws.Range("M35").Value = -10210
$ws.Range("H36").Value = 12331.667
$ws.Range("I36").Value = 3750
$ws.Range("K36").Value = 3750
$ws.Range("M36").Value = -3500
$ws.Range("H37").Value = 25247.5
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("H113").Value = 756.5
$ws.Range("I113").Value = 507.8
$ws.Range("K113").Value = 1523.4
$ws.Range("M113").Value = 646.5999999999999
$ws.Range("H136").Value = 3894.9048
$ws.Range("I136").Value = 1346.625
$ws.Range("K136").Value = 4039.875
$ws.Range("M136").Value = -1489.875
